$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 44: Institution corrected + birth year corrected
$ws.Range("C44").Value = "École Supérieure Privée d'Ingénierie et de Technologies"
$ws.Range("G44").Value = 1990

# Row 49: country name corrected to French form
$ws.Range("D49").Value = "Arabie Saoudite"

# Row 52: country name corrected to French form
$ws.Range("D52").Value = "Suède"

# New row 54: additional researcher entry
$ws.Range("A54").Value = "Ali"
$ws.Range("B54").Value = "Saad"
$ws.Range("C54").Value = "Aarhus University"
$ws.Range("D54").Value = "Danemark"
$ws.Range("E54").Value = "9cdCQAoAAAAJ"
$ws.Range("F54").Value = "M"
$ws.Range("G54").Value = 1988
$ws.Range("H54").Value = "Chimie et Sciences des Matériaux"

# Match formatting used by the other rows in column F (Genre)
$ws.Range("F53").Copy() | Out-Null
$ws.Range("F54").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Update selection/view to match the saved workbook state
$ws.Range("E55").Select() | Out-Null
